$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Corr/total marks update on the concise marksheet
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
